$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$text = $wsHoja1.Range("A1").Value2
$text = $text -replace [regex]::Escape("✅ 1000 Bs = 3.25 = 12426.54 pesos"), "✅ 1000 Bs = 3.28 = 12524.75 pesos"
$text = $text -replace [regex]::Escape("✅ 12426.54 pesos = 3.23 = 969.47 Bs"), "✅ 12524.75 pesos = 3.27 = 963.2 Bs"
$wsHoja1.Range("A1").Value = $text

# --- Sheet "tasas": update the rate cells ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 305
$wsTasas.Range("O10").Value = 3820.05
$wsTasas.Range("N12").Value = 3835.98
$wsTasas.Range("O12").Value = 295
